$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 13.098
$ws.Range("E6").Value = 12.748
$ws.Range("D7").Value = -7.27
$ws.Range("B8").Value = 5.792
$ws.Range("E9").Value = 12.499
$ws.Range("B10").Value = 7.491000000000001
$ws.Range("E10").Value = 12.581
$ws.Range("B12").Value = 5.935
$ws.Range("C13").Value = -12.729
$ws.Range("B18").Value = 6.364
$ws.Range("D20").Value = -8.222

$wb.Save()
